$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new fluid-inclusion measurement ("K21-59-FI#1.txt") needs to be
# inserted as the new first data row (row 2), pushing the existing
# data rows (currently rows 2-9) down to rows 3-10. Copy the existing
# rows downward one at a time (bottom-most first so nothing gets
# clobbered) instead of using Rows.Insert, which would drag the bold
# header-row border formatting down onto the new row.
for ($r = 9; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($dest, 5).Value = $ws.Cells.Item($r, 5).Value()
}

# Populate the newly freed-up row 2 with the new fluid-inclusion entry.
$ws.Cells.Item(2, 1).Value = "K21-59-FI#1.txt"
$ws.Cells.Item(2, 2).Value = 1151.271034507507
$ws.Cells.Item(2, 3).Value = 49.79430471549424
$ws.Cells.Item(2, 4).Value = 26.99710352365537
$ws.Cells.Item(2, 5).Value = "Spline"
